# Update the two-digit multiplication answers in the table.
# Each cell's "A×B=C" text is replaced with a new equation+answer.
#
# NOTE ON ORDERING: the new text for "14×79=1106" is "25×40=1000",
# which is also the *old* text that must become "77×61=4697".
# To avoid the first replacement's output being clobbered (or the second
# replacement accidentally matching text just inserted by the first),
# the "25×40=1000" -> "77×61=4697" replacement runs before the
# "14×79=1106" -> "25×40=1000" replacement.

$d = $word.ActiveDocument

$d.Content.Find.Execute("85×98=8330", $true, $false, $false, $false, $false, $true, 1, $false, "41×34=1394", 2) | Out-Null
$d.Content.Find.Execute("68×48=3264", $true, $false, $false, $false, $false, $true, 1, $false, "39×67=2613", 2) | Out-Null
$d.Content.Find.Execute("43×80=3440", $true, $false, $false, $false, $false, $true, 1, $false, "48×72=3456", 2) | Out-Null
$d.Content.Find.Execute("92×17=1564", $true, $false, $false, $false, $false, $true, 1, $false, "30×44=1320", 2) | Out-Null
$d.Content.Find.Execute("84×71=5964", $true, $false, $false, $false, $false, $true, 1, $false, "90×41=3690", 2) | Out-Null
$d.Content.Find.Execute("90×56=5040", $true, $false, $false, $false, $false, $true, 1, $false, "30×59=1770", 2) | Out-Null
$d.Content.Find.Execute("68×28=1904", $true, $false, $false, $false, $false, $true, 1, $false, "64×26=1664", 2) | Out-Null
$d.Content.Find.Execute("95×31=2945", $true, $false, $false, $false, $false, $true, 1, $false, "16×84=1344", 2) | Out-Null
$d.Content.Find.Execute("49×62=3038", $true, $false, $false, $false, $false, $true, 1, $false, "54×67=3618", 2) | Out-Null
$d.Content.Find.Execute("94×92=8648", $true, $false, $false, $false, $false, $true, 1, $false, "79×31=2449", 2) | Out-Null
$d.Content.Find.Execute("25×40=1000", $true, $false, $false, $false, $false, $true, 1, $false, "77×61=4697", 2) | Out-Null
$d.Content.Find.Execute("14×79=1106", $true, $false, $false, $false, $false, $true, 1, $false, "25×40=1000", 2) | Out-Null
$d.Content.Find.Execute("56×78=4368", $true, $false, $false, $false, $false, $true, 1, $false, "32×98=3136", 2) | Out-Null
$d.Content.Find.Execute("31×11=341", $true, $false, $false, $false, $false, $true, 1, $false, "13×72=936", 2) | Out-Null
$d.Content.Find.Execute("85×46=3910", $true, $false, $false, $false, $false, $true, 1, $false, "93×32=2976", 2) | Out-Null
$d.Content.Find.Execute("11×95=1045", $true, $false, $false, $false, $false, $true, 1, $false, "42×88=3696", 2) | Out-Null
$d.Content.Find.Execute("67×26=1742", $true, $false, $false, $false, $false, $true, 1, $false, "14×35=490", 2) | Out-Null
$d.Content.Find.Execute("30×69=2070", $true, $false, $false, $false, $false, $true, 1, $false, "58×43=2494", 2) | Out-Null
$d.Content.Find.Execute("21×19=399", $true, $false, $false, $false, $false, $true, 1, $false, "84×38=3192", 2) | Out-Null
$d.Content.Find.Execute("93×34=3162", $true, $false, $false, $false, $false, $true, 1, $false, "48×56=2688", 2) | Out-Null
$d.Content.Find.Execute("37×46=1702", $true, $false, $false, $false, $false, $true, 1, $false, "54×41=2214", 2) | Out-Null
$d.Content.Find.Execute("99×54=5346", $true, $false, $false, $false, $false, $true, 1, $false, "52×89=4628", 2) | Out-Null
$d.Content.Find.Execute("89×97=8633", $true, $false, $false, $false, $false, $true, 1, $false, "12×74=888", 2) | Out-Null
$d.Content.Find.Execute("67×81=5427", $true, $false, $false, $false, $false, $true, 1, $false, "64×67=4288", 2) | Out-Null
$d.Content.Find.Execute("50×53=2650", $true, $false, $false, $false, $false, $true, 1, $false, "50×68=3400", 2) | Out-Null
